# "Ejercicio 1, me faltaba pushear"
#
# Two measurements (I8 and I9, the Vout readings for 7000 Hz and 8000 Hz)
# had been entered swapped. This fixes the swap and annotates the rows
# that were touched with a small audit trail in columns L:N:
#   L  -> tag ("Inter" on the row whose value moved IN, "cambiados" on the
#         row that changed the measurement that was wrong)
#   M  -> the original (pre-swap) value, preserved for reference
#   N  -> "original" label

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Write the new annotation text first and in the order the strings should
# end up in the shared-string table: "Inter", "cambiados", "original".
$ws.Range("L8").Value = "Inter"
$ws.Range("L9").Value = "cambiados"
$ws.Range("N8").Value = "original"
$ws.Range("N9").Value = "original"

# Preserve the original (pre-swap) readings in column M.
$ws.Range("M8").Value = 2.2599999999999998
$ws.Range("M9").Value = 2.6

# Correct the swapped Vout measurements in column I; J8/J9 (=20*LOG10(I/H))
# recalculate automatically.
$ws.Range("I8").Value = 2.6
$ws.Range("I9").Value = 2.2599999999999998

# Match the author's final view state: scrolled so column B is left-most,
# with N9 as the active/selected cell.
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("N9").Select()
